$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D (Price) updates: these are numeric-looking text values in the
# source data (stored as text, e.g. "242.97"), so we force the cells to Text
# format before assigning, which keeps the exact string (incl. trailing zeros,
# long decimals) instead of Excel auto-converting them to a Number.
$priceCells = @(
    @{Cell="D2"; Value="243.52"},
    @{Cell="D3"; Value="22.94"},
    @{Cell="D4"; Value="5.406"},
    @{Cell="D5"; Value="0.05915"},
    @{Cell="D6"; Value="3.392"},
    @{Cell="D7"; Value="6.432"},
    @{Cell="D8"; Value="0.8082"},
    @{Cell="D9"; Value="0.9114"},
    @{Cell="D10"; Value="0.1421"},
    @{Cell="D11"; Value="0.07431"},
    @{Cell="D12"; Value="0.03353"},
    @{Cell="D13"; Value="0.03079"},
    @{Cell="D14"; Value="0.09319"},
    @{Cell="D15"; Value="3.945"},
    @{Cell="D16"; Value="0.001587"},
    @{Cell="D17"; Value="0.04807"},
    @{Cell="D18"; Value="0.0005946"},
    @{Cell="D19"; Value="0.005512"},
    @{Cell="D20"; Value="0.004302"},
    @{Cell="D21"; Value="0.0009866"},
    @{Cell="D22"; Value="0.00007521"},
    @{Cell="D23"; Value="3.661"},
    @{Cell="D24"; Value="2.185"},
    @{Cell="D25"; Value="0.3248"},
    @{Cell="D26"; Value="0.1348"},
    @{Cell="D27"; Value="0.0002448"},
    @{Cell="D40"; Value="0.03895"},
    @{Cell="D41"; Value="0.006205"},
    @{Cell="D42"; Value="0.1066"},
    @{Cell="D43"; Value="0.002746"},
    @{Cell="D44"; Value="0.006514"},
    @{Cell="D45"; Value="0.00005175"},
    @{Cell="D46"; Value="0.00000000751"},
    @{Cell="D47"; Value="0.0005805"},
    @{Cell="D48"; Value="1.050"},
    @{Cell="D49"; Value="0.002319"},
    @{Cell="D50"; Value="0.00002102"},
    @{Cell="D51"; Value="0.0002002"}
)
foreach ($item in $priceCells) {
    $rng = $ws.Range($item.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $item.Value
}

# --- Columns B/C/E (text) updates: plain text, no special handling needed.
$textCells = @(
    @{Cell="B18"; Value="One"},
    @{Cell="C18"; Value="https://coinranking.com/coin/6Lga5NiXX3rT+one-one"},
    @{Cell="E18"; Value="17OneONE"},
    @{Cell="B19"; Value="TigerCash"},
    @{Cell="C19"; Value="https://coinranking.com/coin/6hIn06L2+tigercash-tch"},
    @{Cell="E19"; Value="18TigerCashTCH"},
    @{Cell="B20"; Value="HotbitToken"},
    @{Cell="C20"; Value="https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"},
    @{Cell="E20"; Value="19HotbitTokenHTB"},
    @{Cell="B21"; Value="BitKan"},
    @{Cell="C21"; Value="https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"},
    @{Cell="E21"; Value="20BitKanKAN"},
    @{Cell="B22"; Value="NitroEx"},
    @{Cell="C22"; Value="https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"},
    @{Cell="E22"; Value="21NitroExNTX"},
    @{Cell="B23"; Value="LEO"},
    @{Cell="C23"; Value="https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"},
    @{Cell="E23"; Value="22LEOLEO"},
    @{Cell="B24"; Value="BTSEToken"},
    @{Cell="C24"; Value="https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"},
    @{Cell="E24"; Value="23BTSETokenBTSE"},
    @{Cell="E48"; Value="47CoinbaseStockTokenCOINBestin24h"}
)
foreach ($item in $textCells) {
    $ws.Range($item.Cell).Value = $item.Value
}
